$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.969.06'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').Value = '3.416.01'
$ws.Range('E3').Value = '  -0.83%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '410.67'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.32'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.628'
$ws.Range('E7').Value = '  +0.57%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.727'
$ws.Range('E9').Value = '  -1.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.139'
$ws.Range('E10').Value = '  -1.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '43.13'
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '9.15'
$ws.Range('E12').Value = '  +2.05%  '
$ws.Range('D13').Value = '3.962.89'
$ws.Range('E13').Value = '  -0.40%  '
$ws.Range('E14').Value = '  +0.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000213'
$ws.Range('E15').Value = '  +3.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.97'
$ws.Range('E16').Value = '  -1.78%  '
$ws.Range('D17').Value = '3.407.74'
$ws.Range('E17').Value = '  -2.11%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.53'
$ws.Range('E18').Value = '  +0.68%  '
$ws.Range('B19').Value = 'Polygon'
$ws.Range('C19').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.10'
$ws.Range('E19').Value = '  +1.98%  '
$ws.Range('D20').Value = '61.910.59'
$ws.Range('E20').Value = '  -0.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '488.60'
$ws.Range('E21').Value = '  +19.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '91.71'
$ws.Range('E22').Value = '  +1.78%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.28'
$ws.Range('E23').Value = '  +2.52%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.41'
$ws.Range('E24').Value = '  -0.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.33'
$ws.Range('E25').Value = '  +3.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '34.16'
$ws.Range('E26').Value = '  +2.43%  '
$ws.Range('E27').Value = '  +4.85%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '4.81'
$ws.Range('E28').Value = '  +0.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.71'
$ws.Range('E29').Value = '  +1.23%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.68'
$ws.Range('E30').Value = '  -3.59%  '
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '11.99'
$ws.Range('E31').Value = '  +0.48%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.167'
$ws.Range('E32').Value = '  -2.67%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.113'
$ws.Range('E33').Value = '  -4.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '41.69'
$ws.Range('E34').Value = '  -4.76%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '58.46'
$ws.Range('E36').Value = '  +8.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0491'
$ws.Range('E37').Value = '  -2.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.998'
$ws.Range('E38').Value = '  -0.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '151.17'
$ws.Range('E39').Value = '  +6.97%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.41'
$ws.Range('E40').Value = '  +0.11%  '
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.136'
$ws.Range('E41').Value = '  +3.09%  '
$ws.Range('E42').Value = '  +2.84%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.12'
$ws.Range('E43').Value = '  +6.16%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.94'
$ws.Range('E44').Value = '  +0.46%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.63'
$ws.Range('E45').Value = '  +10.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.22'
$ws.Range('E46').Value = '  +4.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.34'
$ws.Range('E47').Value = '  +20.94%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '16.50'
$ws.Range('E48').Value = '  -1.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.56'
$ws.Range('E49').Value = '  +2.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '117.39'
$ws.Range('E50').Value = '  +21.13%  '
$ws.Range('E51').Value = '  +14.47%  '
